$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue $ws 'D2' '63.914.27'
Set-TextValue $ws 'E2' '  +0.89%  '
Set-TextValue $ws 'D3' '2.641.89'
Set-TextValue $ws 'E3' '  -1.31%  '
Set-TextValue $ws 'D5' '608.05'
Set-TextValue $ws 'E5' '  -1.12%  '
Set-TextValue $ws 'D6' '147.34'
Set-TextValue $ws 'E6' '  +2.48%  '
Set-TextValue $ws 'D8' '0.590'
Set-TextValue $ws 'E8' '  +0.40%  '
Set-TextValue $ws 'E9' '  +2.11%  '
Set-TextValue $ws 'D10' '0.385'
Set-TextValue $ws 'E10' '  +6.32%  '
Set-TextValue $ws 'D11' '5.60'
Set-TextValue $ws 'E11' '  -0.37%  '
Set-TextValue $ws 'E12' '  -1.03%  '
Set-TextValue $ws 'D13' '27.49'
Set-TextValue $ws 'E13' '  +0.27%  '
Set-TextValue $ws 'D14' '3.114.68'
Set-TextValue $ws 'E14' '  -1.36%  '
Set-TextValue $ws 'D15' '63.724.68'
Set-TextValue $ws 'E15' '  +0.78%  '
Set-TextValue $ws 'E16' '  +1.64%  '
Set-TextValue $ws 'D17' '2.626.39'
Set-TextValue $ws 'E17' '  -2.32%  '
Set-TextValue $ws 'D18' '11.78'
Set-TextValue $ws 'E18' '  +2.87%  '
Set-TextValue $ws 'D19' '4.57'
Set-TextValue $ws 'E19' '  +3.65%  '
Set-TextValue $ws 'D20' '346.97'
Set-TextValue $ws 'E20' '  +1.28%  '
Set-TextValue $ws 'E21' '  +0.77%  '
Set-TextValue $ws 'E22' '  -0.02%  '
Set-TextValue $ws 'E23' '  -1.29%  '
Set-TextValue $ws 'D24' '66.32'
Set-TextValue $ws 'E24' '  -1.39%  '
Set-TextValue $ws 'E25' '  +8.37%  '
Set-TextValue $ws 'B26' 'InternetComputer(DFINITY)'
Set-TextValue $ws 'C26' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D26' '9.29'
Set-TextValue $ws 'E26' '  +7.20%  '
Set-TextValue $ws 'B27' 'Fetch.AI'
Set-TextValue $ws 'C27' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D27' '1.69'
Set-TextValue $ws 'E27' '  +2.19%  '
Set-TextValue $ws 'D28' '563.17'
Set-TextValue $ws 'E28' '  +4.44%  '
Set-TextValue $ws 'D29' '8.13'
Set-TextValue $ws 'E29' '  +2.68%  '
Set-TextValue $ws 'E30' '  -0.05%  '
Set-TextValue $ws 'E31' '  -2.40%  '
Set-TextValue $ws 'E32' '  -1.28%  '
Set-TextValue $ws 'D33' '0.0₃0855'
Set-TextValue $ws 'E33' '  +5.86%  '
Set-TextValue $ws 'E34' '  -1.42%  '
Set-TextValue $ws 'E35' '  +2.01%  '
Set-TextValue $ws 'D36' '168.95'
Set-TextValue $ws 'E36' '  -1.90%  '
Set-TextValue $ws 'B37' 'PolygonEcosystemToken'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws 'D37' '0.406'
Set-TextValue $ws 'E37' '  -0.07%  '
Set-TextValue $ws 'B38' 'FirstDigitalUSD'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D38' '0.999'
Set-TextValue $ws 'E38' '  -0.03%  '
Set-TextValue $ws 'E39' '  +4.48%  '
Set-TextValue $ws 'D40' '19.19'
Set-TextValue $ws 'E40' '  -0.22%  '
Set-TextValue $ws 'E41' '  +0.02%  '
Set-TextValue $ws 'D42' '165.11'
Set-TextValue $ws 'E42' '  -6.85%  '
Set-TextValue $ws 'D43' '40.06'
Set-TextValue $ws 'E43' '  -0.12%  '
Set-TextValue $ws 'D44' '3.80'
Set-TextValue $ws 'E44' '  +1.21%  '
Set-TextValue $ws 'D45' '21.99'
Set-TextValue $ws 'E45' '  -1.40%  '
Set-TextValue $ws 'D46' '0.0569'
Set-TextValue $ws 'E46' '  -0.48%  '
Set-TextValue $ws 'D47' '0.627'
Set-TextValue $ws 'E47' '  -1.45%  '
Set-TextValue $ws 'B48' 'dogwifhat'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D48' '2.00'
Set-TextValue $ws 'E48' '  +14.17%  '
Set-TextValue $ws 'B49' 'VeChain'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D49' '0.0246'
Set-TextValue $ws 'E49' '  +2.55%  '
Set-TextValue $ws 'D50' '0.0959'
Set-TextValue $ws 'E50' '  -0.54%  '
Set-TextValue $ws 'D51' '18.85'
Set-TextValue $ws 'E51' '  -0.03%  '
